$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: replaced by a Datacenter infrastructure risk -----------------
$ws.Range("C6").Value = "Problemas de infraestrutura/rede na instalação de aplicativos, configuração e/ou manutenção no Datacenter."
$ws.Range("D6").Value = "Identificamos a necessidade de um profissional para solução de possíveis problemas de infraestrutura/redes no Datacenter. "
$ws.Range("J6").Value = "Contratar profissional de infraestrutura/redes em regime CLT."

# --- Row 5: risk about market/requirement changes -----------------------
$ws.Range("J5").Value = "Contratar um especialista da área de desenvolvimento de software que ficará alocado apenas para prospectar as inovações e/ou mudanças no mercado."
$ws.Range("C5").Value = "Realizar correções/atualizações na aplicação principal não previstas devido a possíveis mudanças de mercado."
$ws.Range("D5").Value = "Fatores externos podem vir a prejudicar o desenvolvimento da aplicação, ocasionando correções inesperadas no projeto e que podem comprometer a qualidade final."

# --- Row 6: remaining fields (type D, impact 5) ---------------------------
$ws.Range("E6").Value = "D"
$ws.Range("F6").Value = 5

# --- Row 8: responsible changed -------------------------------------------
$ws.Range("I8").Value = "Gerente de Configuração e Mudança."

# --- Update selection to match the author's last cursor position ---------
$ws.Range("D5").Select()
